$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 846; this shifts the existing rows
# 846-884 down to 848-886 (Excel also carries the trailing two rows of
# data that now land on the newly extended 885-886 rows automatically).
$ws.Rows("846:847").Insert()

# Populate the two newly inserted rows with their data.
$ws.Range("A846").Value = 6
$ws.Range("B846").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C846").Value = "Metropolitana"
$ws.Range("D846").Value = 44753
$ws.Range("E846").Value = 13
$ws.Range("F846").Value = 100112040
$ws.Range("G846").Value = "Cilantro"
$ws.Range("H846").Value = "Sin especificar"
$ws.Range("I846").Value = "Primera"
$ws.Range("J846").Value = 340
$ws.Range("K846").Value = 10000
$ws.Range("L846").Value = 11000
$ws.Range("M846").Value = 10441
$ws.Range("N846").Value = "$/caja 36 atados"
$ws.Range("O846").Value = "Región Metropolitana"
$ws.Range("P846").Value = 290
$ws.Range("Q846").Value = 36
$ws.Range("R846").Value = "Hortaliza"

$ws.Range("A847").Value = 6
$ws.Range("B847").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C847").Value = "Metropolitana"
$ws.Range("D847").Value = 44753
$ws.Range("E847").Value = 13
$ws.Range("F847").Value = 100112040
$ws.Range("G847").Value = "Cilantro"
$ws.Range("H847").Value = "Sin especificar"
$ws.Range("I847").Value = "Primera"
$ws.Range("J847").Value = 290
$ws.Range("K847").Value = 17000
$ws.Range("L847").Value = 18000
$ws.Range("M847").Value = 17448
$ws.Range("N847").Value = "$/docena de atados"
$ws.Range("O847").Value = "Región Metropolitana"
$ws.Range("P847").Value = 5816
$ws.Range("Q847").Value = 3
$ws.Range("R847").Value = "Hortaliza"
